# Fix error in calculating per capita GDP data:
# Each cumulative "NaN count" in column B (rows 2-184) was overstated by 6,
# so subtract 6 from every value in that range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NANCounts")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value()
    $cell.Value = $current - 6
}
